# Update DAMSLTag (column I) and DialogAct (column J) values for the specified rows
# following a re-run of SGNN dialog act annotation after transcript cleanup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(3, 'b', 'Acknowledge (Backchannel)'),
    @(13, 'ba', 'Appreciation'),
    @(17, 'sd', 'Statement-non-opinion'),
    @(18, 'sd', 'Statement-non-opinion'),
    @(19, 'sd', 'Statement-non-opinion'),
    @(21, 'ba', 'Appreciation'),
    @(53, 'sv', 'Statement-opinion'),
    @(54, 'ba', 'Appreciation'),
    @(67, 'ba', 'Appreciation'),
    @(68, 'b', 'Acknowledge (Backchannel)'),
    @(75, 'sv', 'Statement-opinion'),
    @(76, 'sv', 'Statement-opinion'),
    @(84, 'ba', 'Appreciation'),
    @(85, 'b', 'Acknowledge (Backchannel)'),
    @(88, 'aa', 'Agree/Accept'),
    @(93, 'aa', 'Agree/Accept'),
    @(102, '%', 'Uninterpretable'),
    @(108, 'ba', 'Appreciation'),
    @(119, 'sd', 'Statement-non-opinion'),
    @(121, 'sd', 'Statement-non-opinion'),
    @(126, 'sd', 'Statement-non-opinion'),
    @(134, 'sd', 'Statement-non-opinion'),
    @(136, 'aa', 'Agree/Accept'),
    @(137, 'aa', 'Agree/Accept'),
    @(146, 'sd', 'Statement-non-opinion'),
    @(147, 'aa', 'Agree/Accept'),
    @(149, 'sd', 'Statement-non-opinion'),
    @(156, 'sd', 'Statement-non-opinion'),
    @(157, 'aa', 'Agree/Accept'),
    @(159, 'sd', 'Statement-non-opinion'),
    @(165, 'ba', 'Appreciation'),
    @(170, 'ba', 'Appreciation'),
    @(183, 'ba', 'Appreciation'),
    @(187, 'ba', 'Appreciation'),
    @(190, 'ba', 'Appreciation'),
    @(195, 'ba', 'Appreciation'),
    @(201, 'b', 'Acknowledge (Backchannel)'),
    @(206, 'sd', 'Statement-non-opinion'),
    @(209, 'sd', 'Statement-non-opinion'),
    @(211, 'sd', 'Statement-non-opinion'),
    @(213, 'sd', 'Statement-non-opinion'),
    @(214, 'sd', 'Statement-non-opinion'),
    @(224, 'sd', 'Statement-non-opinion'),
    @(230, 'b', 'Acknowledge (Backchannel)'),
    @(233, 'b', 'Acknowledge (Backchannel)'),
    @(238, 'aa', 'Agree/Accept'),
    @(239, 'sv', 'Statement-opinion'),
    @(241, 'ba', 'Appreciation'),
    @(245, 'aa', 'Agree/Accept'),
    @(251, 'ba', 'Appreciation'),
    @(257, 'b', 'Acknowledge (Backchannel)'),
    @(261, 'aa', 'Agree/Accept'),
    @(263, 'sd', 'Statement-non-opinion'),
    @(271, 'sd', 'Statement-non-opinion'),
    @(274, 'ba', 'Appreciation'),
    @(275, 'sd', 'Statement-non-opinion'),
    @(282, 'ba', 'Appreciation'),
    @(284, 'b', 'Acknowledge (Backchannel)'),
    @(285, 'b', 'Acknowledge (Backchannel)'),
    @(287, 'b', 'Acknowledge (Backchannel)'),
    @(289, 'aa', 'Agree/Accept'),
    @(294, 'ba', 'Appreciation'),
    @(312, 'sv', 'Statement-opinion'),
    @(314, 'ba', 'Appreciation'),
    @(319, 'sv', 'Statement-opinion'),
    @(324, 'sd', 'Statement-non-opinion'),
    @(334, 'sd', 'Statement-non-opinion'),
    @(347, 'sd', 'Statement-non-opinion'),
    @(349, 'sd', 'Statement-non-opinion'),
    @(363, 'sv', 'Statement-opinion'),
    @(371, '%', 'Uninterpretable'),
    @(373, '%', 'Uninterpretable'),
    @(385, 'ba', 'Appreciation'),
    @(387, 'sd', 'Statement-non-opinion'),
    @(394, 'b', 'Acknowledge (Backchannel)'),
    @(426, 'sd', 'Statement-non-opinion'),
    @(455, 'sd', 'Statement-non-opinion'),
    @(468, 'ba', 'Appreciation'),
    @(469, 'b', 'Acknowledge (Backchannel)'),
    @(471, 'sd', 'Statement-non-opinion'),
    @(475, 'ba', 'Appreciation'),
    @(480, 'aa', 'Agree/Accept'),
    @(481, 'aa', 'Agree/Accept'),
    @(494, 'sd', 'Statement-non-opinion'),
    @(496, 'sd', 'Statement-non-opinion'),
    @(498, 'sd', 'Statement-non-opinion'),
    @(504, 'sv', 'Statement-opinion'),
    @(511, 'ba', 'Appreciation'),
    @(512, 'sv', 'Statement-opinion'),
    @(515, 'ba', 'Appreciation'),
    @(521, 'sv', 'Statement-opinion'),
    @(535, 'sd', 'Statement-non-opinion'),
    @(538, 'sv', 'Statement-opinion'),
    @(539, 'sv', 'Statement-opinion'),
)

foreach ($u in $updates) {
    $row = $u[0]
    $damslTag = $u[1]
    $dialogAct = $u[2]
    $ws.Cells.Item($row, 9).Value = $damslTag
    $ws.Cells.Item($row, 10).Value = $dialogAct
}
